# Apply the per-cell content updates described by the commit diff.
# Cells in column D hold price strings that are numeric-looking (e.g.
# "6.40", "235.65") but must stay plain text (they were authored as
# inline strings, and some - like thousand-grouped "43.414.12" - are not
# valid numbers at all). A leading apostrophe forces Excel to keep the
# literal text instead of silently coercing it to a Number (which would
# also drop meaningful trailing zeros, e.g. "23.20" -> 23.2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '43.414.12'
$ws.Cells.Item(3, 4).Value = '2.380.38'
$ws.Cells.Item(3, 5).Value = '  +5.84%  '
$ws.Cells.Item(4, 5).Value = '  +0.16%  '
$ws.Cells.Item(5, 4).Value = '''235.65'
$ws.Cells.Item(5, 5).Value = '  +1.40%  '
$ws.Cells.Item(6, 4).Value = '''0.657'
$ws.Cells.Item(6, 5).Value = '  +2.40%  '
$ws.Cells.Item(7, 4).Value = '''72.04'
$ws.Cells.Item(7, 5).Value = '  +14.09%  '
$ws.Cells.Item(8, 5).Value = '  +0.12%  '
$ws.Cells.Item(9, 5).Value = '  +4.70%  '
$ws.Cells.Item(10, 4).Value = '''0.0974'
$ws.Cells.Item(10, 5).Value = '  -0.86%  '
$ws.Cells.Item(11, 4).Value = '''57.28'
$ws.Cells.Item(11, 5).Value = '  -0.06%  '
$ws.Cells.Item(12, 4).Value = '''26.91'
$ws.Cells.Item(12, 5).Value = '  +1.83%  '
$ws.Cells.Item(13, 4).Value = '2.729.57'
$ws.Cells.Item(13, 5).Value = '  +5.71%  '
$ws.Cells.Item(14, 5).Value = '  +0.63%  '
$ws.Cells.Item(15, 4).Value = '''15.95'
$ws.Cells.Item(15, 5).Value = '  +2.85%  '
$ws.Cells.Item(16, 4).Value = '''6.26'
$ws.Cells.Item(16, 5).Value = '  +2.91%  '
$ws.Cells.Item(17, 4).Value = '''0.857'
$ws.Cells.Item(17, 5).Value = '  +3.26%  '
$ws.Cells.Item(18, 4).Value = '2.386.70'
$ws.Cells.Item(18, 5).Value = '  +5.95%  '
$ws.Cells.Item(19, 4).Value = '43.440.83'
$ws.Cells.Item(19, 5).Value = '  -0.94%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0991'
$ws.Cells.Item(20, 5).Value = '  +0.97%  '
$ws.Cells.Item(21, 4).Value = '''6.40'
$ws.Cells.Item(21, 5).Value = '  +5.44%  '
$ws.Cells.Item(22, 4).Value = '''74.62'
$ws.Cells.Item(23, 4).Value = '''251.84'
$ws.Cells.Item(23, 5).Value = '  +2.61%  '
$ws.Cells.Item(24, 5).Value = '  +18.54%  '
$ws.Cells.Item(25, 5).Value = '  -0.02%  '
$ws.Cells.Item(26, 4).Value = '''2.48'
$ws.Cells.Item(26, 5).Value = '  +1.92%  '
$ws.Cells.Item(27, 2).Value = 'Toncoin'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(27, 4).Value = '''2.28'
$ws.Cells.Item(27, 5).Value = '  +2.53%  '
$ws.Cells.Item(28, 2).Value = 'EthereumClassic'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(28, 4).Value = '''23.20'
$ws.Cells.Item(28, 5).Value = '  +10.40%  '
$ws.Cells.Item(29, 2).Value = 'Cosmos'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(29, 4).Value = '''10.02'
$ws.Cells.Item(29, 5).Value = '  +2.48%  '
$ws.Cells.Item(30, 4).Value = '''174.74'
$ws.Cells.Item(30, 5).Value = '  +0.90%  '
$ws.Cells.Item(31, 5).Value = '  +8.99%  '
$ws.Cells.Item(32, 4).Value = '''0.127'
$ws.Cells.Item(32, 5).Value = '  -8.41%  '
$ws.Cells.Item(33, 5).Value = '  +2.20%  '
$ws.Cells.Item(34, 5).Value = '  +4.23%  '
$ws.Cells.Item(35, 5).Value = '  +1.52%  '
$ws.Cells.Item(36, 4).Value = '''5.08'
$ws.Cells.Item(36, 5).Value = '  +3.01%  '
$ws.Cells.Item(37, 5).Value = '  +3.71%  '
$ws.Cells.Item(38, 4).Value = '''2.47'
$ws.Cells.Item(38, 5).Value = '  +8.35%  '
$ws.Cells.Item(39, 5).Value = '  +0.17%  '
$ws.Cells.Item(40, 4).Value = '''0.0256'
$ws.Cells.Item(40, 5).Value = '  +0.85%  '
$ws.Cells.Item(41, 2).Value = 'BinanceUSD'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(41, 4).Value = '''1.00'
$ws.Cells.Item(41, 5).Value = '  -0.13%  '
$ws.Cells.Item(42, 2).Value = 'FraxShare'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(42, 4).Value = '''8.94'
$ws.Cells.Item(42, 5).Value = '  +3.64%  '
$ws.Cells.Item(43, 4).Value = '''18.76'
$ws.Cells.Item(43, 5).Value = '  +10.37%  '
$ws.Cells.Item(44, 4).Value = '''1.19'
$ws.Cells.Item(44, 5).Value = '  +10.40%  '
$ws.Cells.Item(45, 4).Value = '''99.92'
$ws.Cells.Item(45, 5).Value = '  +2.13%  '
$ws.Cells.Item(46, 4).Value = '''1.23'
$ws.Cells.Item(46, 5).Value = '  +3.05%  '
$ws.Cells.Item(47, 4).Value = '''4.46'
$ws.Cells.Item(47, 5).Value = '  +2.45%  '
$ws.Cells.Item(48, 4).Value = '''0.0951'
$ws.Cells.Item(48, 5).Value = '  +0.43%  '
$ws.Cells.Item(49, 4).Value = '1.454.37'
$ws.Cells.Item(49, 5).Value = '  +0.84%  '
$ws.Cells.Item(50, 4).Value = '2.603.23'
$ws.Cells.Item(50, 5).Value = '  +5.93%  '
$ws.Cells.Item(51, 2).Value = 'NEARProtocol'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(51, 4).Value = '''2.26'
$ws.Cells.Item(51, 5).Value = '  -2.13%  '
